# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's sheets, refreshing the figures as of the regenerated output.
#
# Sheet 1 = 展览 (Exhibitions)
# Sheet 2 = 演出 (Performances)
# Sheet 3 = 本地生活 (Local life)
# Sheet 4 = 全部类型 (All types)

$wb = $excel.ActiveWorkbook

$wsExhibitions  = $wb.Worksheets.Item(1)
$wsPerformances = $wb.Worksheets.Item(2)
$wsLocalLife    = $wb.Worksheets.Item(3)
$wsAllTypes     = $wb.Worksheets.Item(4)

# 展览 (Exhibitions)
$wsExhibitions.Range("F5").Value = 1917
$wsExhibitions.Range("F6").Value = 1278
$wsExhibitions.Range("F7").Value = 1597
$wsExhibitions.Range("F15").Value = 1022
$wsExhibitions.Range("F16").Value = 41
$wsExhibitions.Range("F27").Value = 98
$wsExhibitions.Range("F28").Value = 27
$wsExhibitions.Range("F34").Value = 22

# 演出 (Performances)
$wsPerformances.Range("F20").Value = 178
$wsPerformances.Range("F43").Value = 96

# 本地生活 (Local life)
$wsLocalLife.Range("F9").Value = 3072
$wsLocalLife.Range("F10").Value = 596

# 全部类型 (All types)
$wsAllTypes.Range("F8").Value = 3072
$wsAllTypes.Range("F9").Value = 596
$wsAllTypes.Range("F11").Value = 1917
$wsAllTypes.Range("F12").Value = 1278
$wsAllTypes.Range("F19").Value = 1022
$wsAllTypes.Range("F20").Value = 41
$wsAllTypes.Range("F25").Value = 178
$wsAllTypes.Range("F38").Value = 98
$wsAllTypes.Range("F40").Value = 27
$wsAllTypes.Range("F44").Value = 96
$wsAllTypes.Range("F50").Value = 22
